# Text updates as supplied by PM&C.
#
# This script updates the explanatory notes on the "Description" sheet of
# the Housing Remote Indigenous workbook, and makes the Description sheet
# the active / selected sheet (matching the author's final view state).

$wb = $excel.ActiveWorkbook

$wsData = $wb.Worksheets.Item("Data")
$wsDesc = $wb.Worksheets.Item("Description")

# --- Update the three footnote cells on the Description sheet -------------
# B7: replaced note about Victoria's exclusion with a note about NPARIH
#     notional targets.
# B8: "...Remote Housing Strategy." -> "...NPRH."
# B9: previously rich-text ("...Remote Housing Strategy" with no trailing
#     period, split across 3 differently-formatted runs) -> plain text
#     ending in "...NPRH."
$wsDesc.Range("B7").Value = "NPARIH notional targets extend to 2014 for refurbishments and to 2018 for new build houses."
$wsDesc.Range("B8").Value = "Victoria and Tasmania exited the NPARIH in 2014 and are not part of the NPRH."
$wsDesc.Range("B9").Value = "New South Wales exited the NPARIH in 2015 and is not part of the NPRH."

# Unify the formatting of the three footnote cells (B7:B9) so they all use
# the same plain Arial 12pt black font as the other description paragraphs
# (B5/B6), removing the old mixed Calibri/Arial rich-text runs from B9.
$noteFont = $wsDesc.Range("B7:B9").Font
$noteFont.Family = 2
$noteFont.Name = "Arial"
$noteFont.Size = 12
$noteFont.Color = 0

# Row heights shrink now that the notes are shorter / no longer wrap as much.
$wsDesc.Rows.Item(1).RowHeight = 15
$wsDesc.Rows.Item(5).RowHeight = 26.95
$wsDesc.Rows.Item(7).RowHeight = 26.95
$wsDesc.Rows.Item(8).RowHeight = 15
$wsDesc.Rows.Item(9).RowHeight = 15

# Column B widens now that the footnotes are plain uniform text.
$wsDesc.Columns("B").ColumnWidth = 90.35

# --- Switch the active sheet from Data to Description ---------------------
$wsData.Activate()
$wsData.Range("A1").Select()

$wsDesc.Activate()
$wsDesc.Range("B1:B15").Select()
